$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.867.24'
$ws.Range("E2").Value = '  -0.41%  '
$ws.Range("D3").Value = '3.157.62'
$ws.Range("E3").Value = '  +2.43%  '
$ws.Range("E4").Value = '  -0.40%  '
$ws.Range("D5").Value = '''215.90'
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("D6").Value = '''627.11'
$ws.Range("E6").Value = '  +1.53%  '
$ws.Range("D7").Value = '''1.18'
$ws.Range("E7").Value = '  +35.00%  '
$ws.Range("D8").Value = '''0.369'
$ws.Range("E8").Value = '  -1.90%  '
$ws.Range("D9").Value = '''0.999'
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = '3.154.35'
$ws.Range("E10").Value = '  +2.33%  '
$ws.Range("D11").Value = '''0.763'
$ws.Range("E11").Value = '  +12.88%  '
$ws.Range("E12").Value = '  +7.09%  '
$ws.Range("D13").Value = '''5.76'
$ws.Range("E13").Value = '  +7.15%  '
$ws.Range("D14").Value = '''0.0000246'
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D15").Value = '''35.15'
$ws.Range("E15").Value = '  +6.94%  '
$ws.Range("D16").Value = '90.607.71'
$ws.Range("E16").Value = '  -0.31%  '
$ws.Range("D17").Value = '3.742.55'
$ws.Range("E17").Value = '  +2.93%  '
$ws.Range("D18").Value = '3.147.22'
$ws.Range("E18").Value = '  +2.03%  '
$ws.Range("D19").Value = '''3.74'
$ws.Range("E19").Value = '  +7.94%  '
$ws.Range("D20").Value = '''14.65'
$ws.Range("E20").Value = '  +6.73%  '
$ws.Range("D21").Value = '''475.91'
$ws.Range("E21").Value = '  +10.04%  '
$ws.Range("D22").Value = '''0.0000211'
$ws.Range("E22").Value = '  -5.23%  '
$ws.Range("D23").Value = '''9.17'
$ws.Range("E23").Value = '  +8.33%  '
$ws.Range("D24").Value = '''5.29'
$ws.Range("E24").Value = '  +4.17%  '
$ws.Range("B25").Value = 'NEARProtocol'
$ws.Range("C25").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D25").Value = '''5.92'
$ws.Range("E25").Value = '  +7.38%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = '''95.71'
$ws.Range("E26").Value = '  +14.75%  '
$ws.Range("D27").Value = '''12.36'
$ws.Range("E27").Value = '  +4.50%  '
$ws.Range("D28").Value = '3.326.78'
$ws.Range("E28").Value = '  +3.76%  '
$ws.Range("D29").Value = '''0.998'
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("D30").Value = '''9.38'
$ws.Range("E30").Value = '  +8.30%  '
$ws.Range("E31").Value = '  -1.95%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '''0.215'
$ws.Range("E32").Value = '  +55.77%  '
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").Value = '''1.00'
$ws.Range("E33").Value = '  -6.03%  '
$ws.Range("D34").Value = '''27.63'
$ws.Range("E34").Value = '  +20.18%  '
$ws.Range("D35").Value = '''519.48'
$ws.Range("E35").Value = '  +0.95%  '
$ws.Range("E36").Value = '  +5.80%  '
$ws.Range("D37").Value = '''0.146'
$ws.Range("E37").Value = '  +7.20%  '
$ws.Range("D38").Value = '''3.61'
$ws.Range("E38").Value = '  -5.91%  '
$ws.Range("D39").Value = '''6.94'
$ws.Range("E39").Value = '  +1.31%  '
$ws.Range("E40").Value = '  +3.59%  '
$ws.Range("D41").Value = '''0.0918'
$ws.Range("E41").Value = '  +27.80%  '
$ws.Range("E42").Value = '  -0.40%  '
$ws.Range("D43").Value = '''0.425'
$ws.Range("E43").Value = '  +15.83%  '
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("E45").Value = '  +6.07%  '
$ws.Range("D46").Value = '''0.744'
$ws.Range("E46").Value = '  +23.30%  '
$ws.Range("D48").Value = '''4.75'
$ws.Range("E48").Value = '  +11.87%  '
$ws.Range("D49").Value = '''150.98'
$ws.Range("E49").Value = '  +5.36%  '
$ws.Range("D50").Value = '''45.57'
$ws.Range("E50").Value = '  +4.27%  '
$ws.Range("E51").Value = '  +9.62%  '
